# update scripts wuth new tpm
# Recomputed NATMI ligand-receptor TPM-derived metrics for the Efna1-Epha1 sheet.
# Only the numeric measurement columns (G,H,I,J,M,N,O,P,Q,R,S,T) are refreshed;
# identifier/count columns (A-F,K,L) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 23.630375
$ws.Range("H2").Value = 70.89112499999999
$ws.Range("I2").Value = 0.9002398112414131
$ws.Range("J2").Value = 0.9002398112414129
$ws.Range("M2").Value = 2.211928
$ws.Range("N2").Value = 6.635783999999999
$ws.Range("O2").Value = 0.174938892641363
$ws.Range("P2").Value = 0.1749388926413629
$ws.Range("Q2").Value = 52.26868811299999
$ws.Range("R2").Value = 470.4181930169999
$ws.Range("S2").Value = 0.1574869556902424
$ws.Range("T2").Value = 0.1574869556902424

# Row 3
$ws.Range("G3").Value = 23.630375
$ws.Range("H3").Value = 70.89112499999999
$ws.Range("I3").Value = 0.9002398112414131
$ws.Range("J3").Value = 0.9002398112414129
$ws.Range("O3").Value = 0.3056714815357404
$ws.Range("P3").Value = 0.3056714815357404
$ws.Range("Q3").Value = 91.32930414841667
$ws.Range("R3").Value = 821.9637373357499
$ws.Range("S3").Value = 0.275177636839618
$ws.Range("T3").Value = 0.2751776368396179

# Row 4
$ws.Range("G4").Value = 23.630375
$ws.Range("H4").Value = 70.89112499999999
$ws.Range("I4").Value = 0.9002398112414131
$ws.Range("J4").Value = 0.9002398112414129
$ws.Range("O4").Value = 0.5193896258228967
$ws.Range("P4").Value = 0.5193896258228966
$ws.Range("Q4").Value = 155.1845558832916
$ws.Range("R4").Value = 1396.661002949625
$ws.Range("S4").Value = 0.4675752187115527
$ws.Range("T4").Value = 0.4675752187115526

# Row 5
$ws.Range("I5").Value = 0.06214870537054815
$ws.Range("J5").Value = 0.06214870537054815
$ws.Range("M5").Value = 2.211928
$ws.Range("N5").Value = 6.635783999999999
$ws.Range("O5").Value = 0.174938892641363
$ws.Range("P5").Value = 0.1749388926413629
$ws.Range("Q5").Value = 3.608406623519999
$ws.Range("R5").Value = 32.47565961167999
$ws.Range("S5").Value = 0.01087222569661802
$ws.Range("T5").Value = 0.01087222569661802

# Row 6
$ws.Range("I6").Value = 0.06214870537054815
$ws.Range("J6").Value = 0.06214870537054815
$ws.Range("O6").Value = 0.3056714815357404
$ws.Range("P6").Value = 0.3056714815357404
$ws.Range("S6").Value = 0.01899708684614368
$ws.Range("T6").Value = 0.01899708684614368

# Row 7
$ws.Range("I7").Value = 0.06214870537054815
$ws.Range("J7").Value = 0.06214870537054815
$ws.Range("O7").Value = 0.5193896258228967
$ws.Range("P7").Value = 0.5193896258228966
$ws.Range("S7").Value = 0.03227939282778645
$ws.Range("T7").Value = 0.03227939282778645

# Row 8
$ws.Range("I8").Value = 0.03761148338803896
$ws.Range("J8").Value = 0.03761148338803896
$ws.Range("M8").Value = 2.211928
$ws.Range("N8").Value = 6.635783999999999
$ws.Range("O8").Value = 0.174938892641363
$ws.Range("P8").Value = 0.1749388926413629
$ws.Range("Q8").Value = 2.183754673064
$ws.Range("R8").Value = 19.653792057576
$ws.Range("S8").Value = 0.006579711254502556
$ws.Range("T8").Value = 0.006579711254502554

# Row 9
$ws.Range("I9").Value = 0.03761148338803896
$ws.Range("J9").Value = 0.03761148338803896
$ws.Range("O9").Value = 0.3056714815357404
$ws.Range("P9").Value = 0.3056714815357404
$ws.Range("S9").Value = 0.01149675784997876
$ws.Range("T9").Value = 0.01149675784997876

# Row 10
$ws.Range("I10").Value = 0.03761148338803896
$ws.Range("J10").Value = 0.03761148338803896
$ws.Range("O10").Value = 0.5193896258228967
$ws.Range("P10").Value = 0.5193896258228966
$ws.Range("S10").Value = 0.01953501428355765
$ws.Range("T10").Value = 0.01953501428355765
